$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove stray spell-check markers (w:proofErr) by re-typing the spanned
#    text through Find/Replace: Word's engine merges the runs it touches and
#    drops now-orphaned proofErr elements, mirroring what happens when a
#    user retypes/accepts a word in the real app.
# ---------------------------------------------------------------------------

function Retype($needle) {
    $ok = $d.Content.Find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, $needle, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $needle"
    }
}

Retype("a importância dos SIGs no")
Retype("decisão. Fitz em")

Retype("modelo raster, também")
Retype("na matrix, sendo")
Retype("O modelo raster utiliza")
Retype("Para realizar a discretização da imagem")
Retype("Os satélites de observação")
Retype("a limiarização por um tom")

# ---------------------------------------------------------------------------
# 2) Rework the "Arquitetura" / "Modelo de Dados" part of the report:
#    - the two placeholder paragraphs right under the "Arquitetura" heading
#      are replaced by five new paragraphs describing the change-detection
#      pipeline;
#    - the whole "Modelo de Dados" subsection (heading, placeholder text,
#      figure caption and picture) is deleted outright.
# ---------------------------------------------------------------------------

function FindParagraphIndex($text) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        if ($d.Paragraphs($i).Range.Text.Trim() -eq $text) {
            return $i
        }
    }
    return -1
}

# 2a) Delete the "Modelo de Dados" subsection first (it sits after the
#     "Arquitetura" filler, so removing it does not disturb earlier indices).
$idxModelo = FindParagraphIndex("Modelo de Dados")
$idxDetalhes = FindParagraphIndex("Detalhes")
$pModelo = $d.Paragraphs($idxModelo)
$pBeforeDetalhes = $d.Paragraphs($idxDetalhes - 1)
$rModelo = $d.Range($pModelo.Range.Start, $pBeforeDetalhes.Range.End)
$rModelo.Delete()

# 2b) Replace the two placeholder paragraphs below "Arquitetura" with the
#     five new paragraphs. ("Modelo de Dados" is already gone at this point,
#     so "Detalhes" -- now immediately below the filler -- is the boundary.)
$idxArquitetura = FindParagraphIndex("Arquitetura")
$pFirst = $d.Paragraphs($idxArquitetura + 1)
$idxDetalhesNow = FindParagraphIndex("Detalhes")
$pLast = $d.Paragraphs($idxDetalhesNow - 1)
$rFiller = $d.Range($pFirst.Range.Start, $pLast.Range.End)

$newXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/><w:ind w:firstLine="0"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">A partir da imagem matricial é possível fazer operações para extrair dados e informações sobre uma área de específica. Para a detecção de mudança no uso e cobertura da terra, busca-se identificar as possíveis alterações em um </w:t></w:r><w:r><w:t>período</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/><w:ind w:firstLine="0"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">Neste processo, define-se um ponto de interesse e imagens multitemporais para realizar a comparação e gerar uma nova imagem que represente a diferença detectada. Dessa forma, o sistema desenvolvido realiza as operações de identificação </w:t></w:r><w:r><w:t>da série histórica</w:t></w:r><w:r><w:t>, este recebe 3 imagens como entrada. O primeiro parâmetro é a imagem da área em que se busca realizar a análise, a referência. Os outros parâmetros são as cenas em períodos distintos, o qual a imagem de referência deve abranger.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/><w:ind w:firstLine="0"/><w:jc w:val="left"/></w:pPr><w:r><w:t>Após, um algoritmo para tornar as imagens temporais em escala de cinza é aplicado, a fim de tornar a imagem binária, este processo é feito pela ponderação da soma dos valores dos canais do pixel, o RGB (</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>red, green, blue</w:t></w:r><w:r><w:t>). Assim, realizar a médias das intensidades em cada faixa de frequência para transformá-la em monocromática. Ainda, na segmentação definiu-se níveis para alterar a cor do pixel, para cores abaixo de 0.3 torna-se preto e acima branco.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/><w:ind w:firstLine="0"/><w:jc w:val="left"/></w:pPr><w:r><w:t>Em seguida, é realizado a etapa de identificação da área de interesse nas imagens da comparação, em um processo chamado de Casamento de Modelo (T</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>emplate Matching</w:t></w:r><w:r><w:t>). Para esta operação, assume-se que não há distorções nas imagens. O cálculo com base no número de linhas e colunas da imagem matricial retorna as coordenadas no plano cartesiano do deslocamento, que é usado para alinhá-las.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/><w:ind w:firstLine="0"/><w:jc w:val="left"/></w:pPr><w:r><w:t>Por fim, com o alinhamento das imagens com a área de referência, realiza-se a comparação pixel a pixel das duas imagens temporais, destacando as diferenças entre os pixels e gerando uma nova imagem a partir disso. O resultado da comparação é uma imagem binarizada referenciada em um plano.</w:t></w:r></w:p>
'@

$rFiller.InsertXML($newXml)

Write-Output "edit complete"
